$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.0665
$ws.Range("E2").Value = 0.108
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 121.2
$ws.Range("L2").Value = 0.2614886731391586
$ws.Range("M2").Value = 80.44
$ws.Range("N2").Value = 0.04442971554819111
$ws.Range("O2").Value = 0.6636963696369637
$ws.Range("P2").Value = 70.90000000000001
$ws.Range("Q2").Value = 0.0391604529135598
$ws.Range("R2").Value = 0.584983498349835
$ws.Range("S2").Value = 9.539999999999992
$ws.Range("T2").Value = 0.1185977125808055
$ws.Range("U2").Value = 7778.4
$ws.Range("V2").Value = 4.296271748135874
$ws.Range("W2").Value = 0.06298721546616776
$ws.Range("X2").Value = 0.07525485964080228
$ws.Range("Y2").Value = -0.01226764417463452
$ws.Range("Z2").Value = -0.451050992604126
$ws.Range("AA2").Value = -0
$ws.Range("AB2").Value = 0.03513010433165432
$ws.Range("AC2").Value = -0.03513010433165432
$ws.Range("AD2").Value = 3963.2
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 3963.2
$ws.Range("AG2").Value = -3815.2
$ws.Range("AH2").Value = 0.6864229177130783
$ws.Range("AI2").Value = 0.6473172723560636
$ws.Range("AJ2").Value = 1.903127650022447
$ws.Range("AK2").Value = 2.304003864967691
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()

# Row 3
$ws.Range("D3").Value = 0.0665
$ws.Range("E3").Value = 0.108
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 121.2
$ws.Range("L3").Value = 0.2614886731391586
$ws.Range("M3").Value = 80.44
$ws.Range("N3").Value = 0.04442971554819111
$ws.Range("O3").Value = 0.6636963696369637
$ws.Range("P3").Value = 70.90000000000001
$ws.Range("Q3").Value = 0.0391604529135598
$ws.Range("R3").Value = 0.584983498349835
$ws.Range("S3").Value = 9.539999999999992
$ws.Range("T3").Value = 0.1185977125808055
$ws.Range("U3").Value = 7778.4
$ws.Range("V3").Value = 4.296271748135874
$ws.Range("W3").Value = 0.06298721546616776
$ws.Range("X3").Value = 0.07525485964080228
$ws.Range("Y3").Value = -0.01226764417463452
$ws.Range("Z3").Value = -0.451050992604126
$ws.Range("AA3").Value = -0
$ws.Range("AB3").Value = 0.03513010433165432
$ws.Range("AC3").Value = -0.03513010433165432
$ws.Range("AD3").Value = 3963.2
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 3963.2
$ws.Range("AG3").Value = -3815.2
$ws.Range("AH3").Value = 0.6864229177130783
$ws.Range("AI3").Value = 0.6473172723560636
$ws.Range("AJ3").Value = 1.903127650022447
$ws.Range("AK3").Value = 2.304003864967691
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()
